# Update "BOM Proposal and Options.xlsx":
#  - add a new worksheet "Hot Ends" after "Ark1" containing a water-cooled
#    hot-end comparison table (with a hyperlink + currency formatting)
#  - minor view-state tweaks on the "Ark1" sheet

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add the new "Hot Ends" worksheet right after "Ark1"
# ---------------------------------------------------------------------
$ark1 = $wb.Worksheets.Item(1)
$hotEnds = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ark1)
$hotEnds.Name = "Hot Ends"

# ---------------------------------------------------------------------
# 2. Column widths for the new sheet
# (input values pre-compensated for the engine's internal 1/6-character
# rounding so the stored OOXML width lands as close as possible to the
# original authored widths)
# ---------------------------------------------------------------------
$hotEnds.Columns.Item(1).ColumnWidth = 20.66666666666667
$hotEnds.Columns.Item(2).ColumnWidth = 14.5
$hotEnds.Columns.Item(3).ColumnWidth = 31.16666666666667
$hotEnds.Columns.Item(4).ColumnWidth = 10.5
$hotEnds.Columns.Item(5).ColumnWidth = 15.333333333333334
$hotEnds.Columns.Item(6).ColumnWidth = 13.166666666666666
$hotEnds.Columns.Item(7).ColumnWidth = 11.666666666666666

# ---------------------------------------------------------------------
# 3. Cell content
# ---------------------------------------------------------------------
$hotEnds.Range("A1").Value = "Water cooled Hot Ends"

$hotEnds.Range("B2").Value = "Manufacture "
$hotEnds.Range("C2").Value = "Model"
$hotEnds.Range("D2").Value = "Cost US `$"
$hotEnds.Range("E2").Value = "Max T deg C"
$hotEnds.Range("F2").Value = "Silicon sock"
$hotEnds.Range("G2").Value = "URL"

$hotEnds.Range("B4").Value = "TriangleLabs"
$hotEnds.Range("C4").Value = "Dragon Hotend  standard flow"
$hotEnds.Range("D4").Value = 86
$hotEnds.Range("E4").Value = 500
$hotEnds.Range("F4").Value = "yes"
$hotEnds.Range("G4").Value = "here"

$hotEnds.Range("B5").Value = "TriangleLabs"
$hotEnds.Range("C5").Value = "Dragon Hotend  High flow"
$hotEnds.Range("D5").Value = 96
$hotEnds.Range("E5").Value = 500
$hotEnds.Range("F5").Value = "yes"

# ---------------------------------------------------------------------
# 4. Formatting
# ---------------------------------------------------------------------
# Header row + plain centered cells
$hotEnds.Range("B2:G2").HorizontalAlignment = -4108
$hotEnds.Range("E4:F4").HorizontalAlignment = -4108
$hotEnds.Range("E5:F5").HorizontalAlignment = -4108

# Currency columns, centered
$hotEnds.Range("D4:D5").HorizontalAlignment = -4108
$hotEnds.Range("D4:D5").NumberFormat = "$#,##0_);[Red]($#,##0)"

# Hyperlink cell gets centered alignment before the hyperlink style is applied.
# The cell already holds the display text "here"; re-use it (omit
# TextToDisplay) so the existing shared-string text is kept untouched.
$hotEnds.Range("G4").HorizontalAlignment = -4108
$hotEnds.Hyperlinks.Add($hotEnds.Range("G4"), "https://www.triangle-lab.com/products/dragon-hotend")

# ---------------------------------------------------------------------
# 5. View state
# ---------------------------------------------------------------------
$hotEnds.Range("B8").Select()

$wb.Save()
